{"js": "// The edit turns the first occurrence of the quoted research question\n// (paragraph \"1ste versie onderzoeksvraag:\") from\n//   \"... binnen een grotere entiteit zoals Het Kontakt, ...\"\n// into\n//   \"... binnen een grotere entiteit zoaals Het Kontakt, ...\"\n// i.e. an \"a\" is inserted right after \"zo\" / before \"als Het Kontakt\",\n// while the second (otherwise identical) occurrence later in the\n// document, inside the \"Inleiding\" paragraph, must stay untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Find the one paragraph that is exactly the quoted research question\n// (the short \"1ste versie onderzoeksvraag\" version), not the longer\n// \"Inleiding\" paragraph that happens to repeat the same sentence.\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  // The short \"1ste versie onderzoeksvraag\" paragraph IS the quoted\n  // sentence (starts with the opening quote mark); the later\n  // \"Inleiding\" paragraph only ends with the same quoted sentence.\n  if (t.indexOf(\"zoals Het Kontakt\") !== -1 && t.indexOf('\"Wat zijn de implicaties') === 0) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not locate the target paragraph\");\n}\n\n// Search within that single paragraph only, so the identical sentence\n// that reappears later in the document is left alone.\nconst hits = target.search(\"grotere entiteit zoals Het Kontakt\", { matchCase: true });\nhits.load(\"items\");\nawait context.sync();\n\nif (hits.items.length !== 1) {\n  throw new Error(\"Expected exactly one match in the target paragraph, found \" + hits.items.length);\n}\n\nhits.items[0].insertText(\"grotere entiteit zoaals Het Kontakt\", \"Replace\");\nawait context.sync();\n", "ps1": "# The edit turns the first occurrence of the quoted research question\n# (paragraph \"1ste versie onderzoeksvraag:\") from\n#   \"... binnen een grotere entiteit zoals Het Kontakt, ...\"\n# into\n#   \"... binnen een grotere entiteit zoaals Het Kontakt, ...\"\n# i.e. an \"a\" is inserted right after \"zo\" / before \"als Het Kontakt\".\n# The document also contains a second, otherwise identical, copy of the\n# same sentence further down (inside the \"Inleiding\" paragraph) that\n# must be left untouched, so the edit is scoped to a single paragraph\n# found by its distinctive leading quote mark.\n\n$d = $word.ActiveDocument\n\n$targetParagraph = $null\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text\n    if ($t.StartsWith('\"Wat zijn de implicaties') -and $t.Contains('zoals Het Kontakt')) {\n        $targetParagraph = $p\n        break\n    }\n}\n\nif ($targetParagraph -eq $null) {\n    throw \"Could not locate the target paragraph\"\n}\n\n$searchRange = $targetParagraph.Range\n$found = $searchRange.Find.Execute(\n    \"grotere entiteit zoals Het Kontakt\",\n    $false,\n    $true,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    \"grotere entiteit zoaals Het Kontakt\",\n    2\n)\n\nif (-not $found) {\n    throw \"Could not find the text to replace in the target paragraph\"\n}\n"}
